$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.54709974267697
$ws.Range("C2").Value = 6.314146918713713
$ws.Range("D2").Value = 3.942075114654553
$ws.Range("E2").Value = 6.755857619866471
$ws.Range("F2").Value = 30.04037933696119
$ws.Range("G2").Value = 38.91621182352387
$ws.Range("H2").Value = 3.777892660849062
$ws.Range("J2").Value = 13.06054131139027
$ws.Range("K2").Value = 21.31968121917402
$ws.Range("L2").Value = 5.715269590898211
$ws.Range("M2").Value = 11.16061223456801
$ws.Range("N2").Value = 6.705636065475505
$ws.Range("B3").Value = 12.71533305148697
$ws.Range("C3").Value = 6.00126565964232
$ws.Range("D3").Value = 3.766864689019111
$ws.Range("E3").Value = 6.573162568667697
$ws.Range("F3").Value = 29.6144930922619
$ws.Range("G3").Value = 38.29744678025261
$ws.Range("H3").Value = 4.038738876163019
$ws.Range("J3").Value = 13.02792708849563
$ws.Range("K3").Value = 21.25751191497488
$ws.Range("L3").Value = 5.656553904475721
$ws.Range("M3").Value = 10.45894366041986
$ws.Range("N3").Value = 6.484463046517082
$ws.Range("B4").Value = 12.17868700645685
$ws.Range("C4").Value = 5.800616644023948
$ws.Range("D4").Value = 3.655369680183485
$ws.Range("E4").Value = 6.457984621404178
$ws.Range("F4").Value = 29.36318737143544
$ws.Range("G4").Value = 37.9336017990773
$ws.Range("H4").Value = 4.204297181135205
$ws.Range("J4").Value = 13.01282307414971
$ws.Range("K4").Value = 21.22776277200052
$ws.Range("L4").Value = 5.61948842184915
$ws.Range("M4").Value = 10.00371764795738
$ws.Range("N4").Value = 6.344859415062411
$ws.Range("B5").Value = 11.95357802100971
$ws.Range("C5").Value = 5.716737077658332
$ws.Range("D5").Value = 3.60899302774706
$ws.Range("E5").Value = 6.410315621738614
$ws.Range("F5").Value = 29.2633655418075
$ws.Range("G5").Value = 37.78940311690961
$ws.Range("H5").Value = 4.273145968733564
$ws.Range("J5").Value = 13.00787541893621
$ws.Range("K5").Value = 21.21770194010508
$ws.Range("L5").Value = 5.604125995985633
$ws.Range("M5").Value = 9.812075878749024
$ws.Range("N5").Value = 6.287065043526671
$ws.Range("B6").Value = 11.91581455147016
$ws.Range("C6").Value = 5.702682453190258
$ws.Range("D6").Value = 3.601236748059891
$ws.Range("E6").Value = 6.402356636650103
$ws.Range("F6").Value = 29.24694708454945
$ws.Range("G6").Value = 37.76570553167817
$ws.Range("H6").Value = 4.284662470814403
$ws.Range("J6").Value = 13.00712589937913
$ws.Range("K6").Value = 21.21615425182327
$ws.Range("L6").Value = 5.601559413777155
$ws.Range("M6").Value = 9.779883690952133
$ws.Range("N6").Value = 6.277415295873754
$ws.Range("B7").Value = 12.17567695273143
$ws.Range("C7").Value = 5.799493921143735
$ws.Range("D7").Value = 3.654747978722539
$ws.Range("E7").Value = 6.4573446778773
$ws.Range("F7").Value = 29.36183063472124
$ws.Range("G7").Value = 37.93164057215242
$ws.Range("H7").Value = 4.205220069912869
$ws.Range("J7").Value = 13.01275149980785
$ws.Range("K7").Value = 21.22761881263988
$ws.Range("L7").Value = 5.61928228703473
$ws.Range("M7").Value = 10.00115792555421
$ws.Range("N7").Value = 6.34408357138094
$ws.Range("B8").Value = 13.26574176464058
$ws.Range("C8").Value = 6.208055097806784
$ws.Range("D8").Value = 3.882498758460315
$ws.Range("E8").Value = 6.693499867339094
$ws.Range("F8").Value = 29.89143527452506
$ws.Range("G8").Value = 38.69954391274364
$ws.Range("H8").Value = 3.866727609592196
$ws.Range("J8").Value = 13.04825646022982
$ws.Range("K8").Value = 21.29646360591656
$ws.Range("L8").Value = 5.69523295261534
$ws.Range("M8").Value = 10.92373920311925
$ws.Range("N8").Value = 6.630193967241514
$ws.Range("B9").Value = 15.1944753807414
$ws.Range("C9").Value = 6.940457779380573
$ws.Range("D9").Value = 4.296511391470168
$ws.Range("E9").Value = 7.131957356974689
$ws.Range("F9").Value = 31.00955618043126
$ws.Range("G9").Value = 40.33139506818705
$ws.Range("H9").Value = 3.244579920259121
$ws.Range("J9").Value = 13.15827197703336
$ws.Range("K9").Value = 21.50077932253189
$ws.Range("L9").Value = 5.836278771045402
$ws.Range("M9").Value = 12.53968579566027
$ws.Range("N9").Value = 7.159108497454233
$ws.Range("B10").Value = 16.53129321111028
$ws.Range("C10").Value = 7.460973093141031
$ws.Range("D10").Value = 4.562371903879074
$ws.Range("E10").Value = 7.358589822077299
$ws.Range("F10").Value = 31.65327879244567
$ws.Range("G10").Value = 41.26281081736548
$ws.Range("H10").Value = 2.839264103485546
$ws.Range("J10").Value = 13.19841604468296
$ws.Range("K10").Value = 21.55848081411547
$ws.Range("L10").Value = 5.905821799073609
$ws.Range("M10").Value = 13.61342451234055
$ws.Range("N10").Value = 7.454276240908461
$ws.Range("B11").Value = 17.03466291079945
$ws.Range("C11").Value = 7.909410692668758
$ws.Range("D11").Value = 4.533101463183519
$ws.Range("E11").Value = 6.822175492029073
$ws.Range("F11").Value = 30.02270593079822
$ws.Range("G11").Value = 38.75702818976601
$ws.Range("H11").Value = 3.523400830817337
$ws.Range("J11").Value = 12.64544752814321
$ws.Range("K11").Value = 20.41482141769366
$ws.Range("L11").Value = 5.727302326709221
$ws.Range("M11").Value = 14.10695879722952
$ws.Range("N11").Value = 6.978424061314503
$ws.Range("B12").Value = 17.18370427108654
$ws.Range("C12").Value = 8.189266818028493
$ws.Range("D12").Value = 4.450820185862972
$ws.Range("E12").Value = 6.39756071039896
$ws.Range("F12").Value = 28.49981311100047
$ws.Range("G12").Value = 36.42429589392698
$ws.Range("H12").Value = 4.694251594327262
$ws.Range("J12").Value = 12.17144391047034
$ws.Range("K12").Value = 19.45038963609268
$ws.Range("L12").Value = 5.627164230101688
$ws.Range("M12").Value = 14.30753647328019
$ws.Range("N12").Value = 6.527851283741581
$ws.Range("B13").Value = 17.07445895098845
$ws.Range("C13").Value = 8.362569675739152
$ws.Range("D13").Value = 4.320777977223814
$ws.Range("E13").Value = 6.043742318397285
$ws.Range("F13").Value = 26.95621031826871
$ws.Range("G13").Value = 34.06081452199334
$ws.Range("H13").Value = 6.011652916742624
$ws.Range("J13").Value = 11.72995020749395
$ws.Range("K13").Value = 18.56950139799687
$ws.Range("L13").Value = 5.582263050279278
$ws.Range("M13").Value = 14.30115490478224
$ws.Range("N13").Value = 6.066991862293138
$ws.Range("B14").Value = 16.87988185987772
$ws.Range("C14").Value = 8.438298569458652
$ws.Range("D14").Value = 4.206643622258033
$ws.Range("E14").Value = 5.84899576148434
$ws.Range("F14").Value = 25.87157969850933
$ws.Range("G14").Value = 32.39877457910541
$ws.Range("H14").Value = 6.976679369592144
$ws.Range("J14").Value = 11.43969340818208
$ws.Range("K14").Value = 18.00261935478401
$ws.Range("L14").Value = 5.584790392069773
$ws.Range("M14").Value = 14.203364382628
$ws.Range("N14").Value = 5.74395361777802
$ws.Range("B15").Value = 16.77637493785333
$ws.Range("C15").Value = 8.433324113045069
$ws.Range("D15").Value = 4.168681616385773
$ws.Range("E15").Value = 5.806267889283668
$ws.Range("F15").Value = 25.60579795509562
$ws.Range("G15").Value = 31.99220057136209
$ws.Range("H15").Value = 7.207019582353172
$ws.Range("J15").Value = 11.37646816901481
$ws.Range("K15").Value = 17.88231075155671
$ws.Range("L15").Value = 5.587970026686676
$ws.Range("M15").Value = 14.1346034099927
$ws.Range("N15").Value = 5.66145657072157
$ws.Range("B16").Value = 16.23784104417389
$ws.Range("C16").Value = 8.194957561659951
$ws.Range("D16").Value = 4.077629634522183
$ws.Range("E16").Value = 5.770141749272704
$ws.Range("F16").Value = 25.61647062119166
$ws.Range("G16").Value = 32.02039690608884
$ws.Range("H16").Value = 7.059215143317187
$ws.Range("J16").Value = 11.44834329944642
$ws.Range("K16").Value = 18.03841489063852
$ws.Range("L16").Value = 5.562821099129176
$ws.Range("M16").Value = 13.69781194330178
$ws.Range("N16").Value = 5.615723886521653
$ws.Range("B17").Value = 15.9385786626722
$ws.Range("C17").Value = 7.965838328586629
$ws.Range("D17").Value = 4.069912409129359
$ws.Range("E17").Value = 5.835719250074592
$ws.Range("F17").Value = 26.2044593381743
$ws.Range("G17").Value = 32.93259352822059
$ws.Range("H17").Value = 6.4130770664218
$ws.Range("J17").Value = 11.6564399088733
$ws.Range("K17").Value = 18.453809100753
$ws.Range("L17").Value = 5.536756770117194
$ws.Range("M17").Value = 13.40961455579998
$ws.Range("N17").Value = 5.751819852292315
$ws.Range("B18").Value = 15.79205675693214
$ws.Range("C18").Value = 7.723786438530449
$ws.Range("D18").Value = 4.133021234151022
$ws.Range("E18").Value = 6.055349588442717
$ws.Range("F18").Value = 27.36375813755003
$ws.Range("G18").Value = 34.71886864612912
$ws.Range("H18").Value = 5.317936040439876
$ws.Range("J18").Value = 12.01026162547257
$ws.Range("K18").Value = 19.15808004244545
$ws.Range("L18").Value = 5.541563166865877
$ws.Range("M18").Value = 13.22354720453932
$ws.Range("N18").Value = 6.075739979473928
$ws.Range("B19").Value = 15.7802317465559
$ws.Range("C19").Value = 7.496543904045444
$ws.Range("D19").Value = 4.246574070521728
$ws.Range("E19").Value = 6.458453306233318
$ws.Range("F19").Value = 28.91382806154399
$ws.Range("G19").Value = 37.0985281837162
$ws.Range("H19").Value = 4.096395108887572
$ws.Range("J19").Value = 12.46532242789947
$ws.Range("K19").Value = 20.07405504923958
$ws.Range("L19").Value = 5.620933961358085
$ws.Range("M19").Value = 13.1378036647913
$ws.Range("N19").Value = 6.544734592439132
$ws.Range("B20").Value = 16.18442887020872
$ws.Range("C20").Value = 7.328513578491409
$ws.Range("D20").Value = 4.493213741223263
$ws.Range("E20").Value = 7.296258113205539
$ws.Range("F20").Value = 31.47076670744796
$ws.Range("G20").Value = 40.99719227961923
$ws.Range("H20").Value = 2.944206536227163
$ws.Range("J20").Value = 13.18306225335668
$ws.Range("K20").Value = 21.53401825380146
$ws.Range("L20").Value = 5.88644064301378
$ws.Range("M20").Value = 13.33949870848997
$ws.Range("N20").Value = 7.374774060201046
$ws.Range("B21").Value = 17.19339268536471
$ws.Range("C21").Value = 7.675094051510034
$ws.Range("D21").Value = 4.716220850810917
$ws.Range("E21").Value = 7.58862514267191
$ws.Range("F21").Value = 32.32683630973685
$ws.Range("G21").Value = 42.2650677013775
$ws.Range("H21").Value = 2.597774530425203
$ws.Range("J21").Value = 13.32672359068395
$ws.Range("K21").Value = 21.80744520105378
$ws.Range("L21").Value = 5.984240926343857
$ws.Range("M21").Value = 14.12169316801922
$ws.Range("N21").Value = 7.705062759628508
$ws.Range("B22").Value = 17.82026792756953
$ws.Range("C22").Value = 7.906808185138907
$ws.Range("D22").Value = 4.84911170968028
$ws.Range("E22").Value = 7.735726288158483
$ws.Range("F22").Value = 32.77923888732263
$ws.Range("G22").Value = 42.93108427587725
$ws.Range("H22").Value = 2.38867410706199
$ws.Range("J22").Value = 13.39167522654647
$ws.Range("K22").Value = 21.92491151498069
$ws.Range("L22").Value = 6.032259751895056
$ws.Range("M22").Value = 14.6138459251333
$ws.Range("N22").Value = 7.879656538857756
$ws.Range("B23").Value = 17.48819854764229
$ws.Range("C23").Value = 7.783848734544178
$ws.Range("D23").Value = 4.778572483241211
$ws.Range("E23").Value = 7.657496046790967
$ws.Range("F23").Value = 32.53707123651547
$ws.Range("G23").Value = 42.57446587435918
$ws.Range("H23").Value = 2.499945145495821
$ws.Range("J23").Value = 13.35650960784429
$ws.Range("K23").Value = 21.86135052033183
$ws.Range("L23").Value = 6.006698301377481
$ws.Range("M23").Value = 14.3530657196486
$ws.Range("N23").Value = 7.786886214512292
$ws.Range("B24").Value = 16.17518763773635
$ws.Range("C24").Value = 7.302510078851402
$ws.Range("D24").Value = 4.502975283454851
$ws.Range("E24").Value = 7.355083686443075
$ws.Range("F24").Value = 31.63574361808225
$ws.Range("G24").Value = 41.24942626401924
$ws.Range("H24").Value = 2.928818624847668
$ws.Range("J24").Value = 13.23411928999113
$ws.Range("K24").Value = 21.63933876651797
$ws.Range("L24").Value = 5.908384120025908
$ws.Range("M24").Value = 13.32367808225213
$ws.Range("N24").Value = 7.426573075666864
$ws.Range("B25").Value = 14.69540054807386
$ws.Range("C25").Value = 6.749799246263108
$ws.Range("D25").Value = 4.188190363896424
$ws.Range("E25").Value = 7.0160829508854
$ws.Range("F25").Value = 30.69874092736322
$ws.Range("G25").Value = 39.87668074330681
$ws.Range("H25").Value = 3.408605331343533
$ws.Range("J25").Value = 13.12399557801251
$ws.Range("K25").Value = 21.43773613225504
$ws.Range("L25").Value = 5.798944176778493
$ws.Range("M25").Value = 12.12313473832127
$ws.Range("N25").Value = 7.019717813673934
